$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.847.80"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.578.71"
$ws.Range("E3").Value = "  +2.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'302.35"
$ws.Range("E5").Value = "  +1.93%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'97.01"
$ws.Range("E6").Value = "  +3.79%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.79%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = "  +0.46%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'36.59"
$ws.Range("E10").Value = "  +1.29%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0809"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12
$ws.Range("D12").Value = "'7.70"
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("E13").Value = "  +6.84%  "

# Row 14
$ws.Range("D14").Value = "2.556.31"
$ws.Range("E14").Value = "  +1.84%  "

# Row 15
$ws.Range("D15").Value = "'0.882"
$ws.Range("E15").Value = "  +2.04%  "

# Row 16
$ws.Range("D16").Value = "'14.35"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17
$ws.Range("D17").Value = "42.875.18"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18
$ws.Range("D18").Value = "'12.91"
$ws.Range("E18").Value = "  +5.47%  "

# Row 19
$ws.Range("E19").Value = "  +2.73%  "

# Row 20
$ws.Range("E20").Value = "  +0.57%  "

# Row 21
$ws.Range("D21").Value = "'71.95"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("D22").Value = "'254.35"
$ws.Range("E22").Value = "  -1.98%  "

# Row 23
$ws.Range("D23").Value = "'2.96"
$ws.Range("E23").Value = "  +2.22%  "

# Row 24
$ws.Range("E24").Value = "  -1.89%  "

# Row 25
$ws.Range("D25").Value = "'28.60"
$ws.Range("E25").Value = "  -1.49%  "

# Row 26
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("D27").Value = "'10.21"
$ws.Range("E27").Value = "  +2.47%  "

# Row 28
$ws.Range("D28").Value = "'39.00"
$ws.Range("E28").Value = "  +6.08%  "

# Row 29
$ws.Range("E29").Value = "  -0.72%  "

# Row 30
$ws.Range("E30").Value = "  +0.91%  "

# Row 31
$ws.Range("D31").Value = "'155.30"
$ws.Range("E31").Value = "  +2.53%  "

# Row 32
$ws.Range("E32").Value = "  -1.06%  "

# Row 33
$ws.Range("E33").Value = "  +0.64%  "

# Row 34
$ws.Range("D34").Value = "'0.0812"
$ws.Range("E34").Value = "  +1.53%  "

# Row 35
$ws.Range("E35").Value = "  -3.66%  "

# Row 36
$ws.Range("D36").Value = "'18.46"
$ws.Range("E36").Value = "  +12.29%  "

# Row 37
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("E38").Value = "  +0.71%  "

# Row 39
$ws.Range("D39").Value = "'23.37"
$ws.Range("E39").Value = "  -1.75%  "

# Row 40
$ws.Range("E40").Value = "  +30.44%  "

# Row 41
$ws.Range("E41").Value = "  +1.60%  "

# Row 42
$ws.Range("E42").Value = "  -1.26%  "

# Row 43
$ws.Range("E43").Value = "  +0.54%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.069.40"
$ws.Range("E44").Value = "  +2.21%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  +3.79%  "

# Row 47 - BitcoinSV
$ws.Range("D47").Value = "'85.23"
$ws.Range("E47").Value = "  -0.46%  "

# Row 48 - was RocketPoolETH, now ordi
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'76.16"
$ws.Range("E48").Value = "  +11.23%  "

# Row 49 - was ordi, now RocketPoolETH
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'2.831.51"
$ws.Range("E49").Value = "  +2.55%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'106.03"
$ws.Range("E50").Value = "  +2.93%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  +2.48%  "

